$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "1.00", "0.519").
# Excel auto-converts such strings to numbers on assignment, which would lose the
# original text formatting (trailing zeros, decimal count, etc). Temporarily mark
# them as Text format so the assignment keeps the literal text, then reset the
# cell style back to "Normal" (the workbook's default, unstyled look) so the
# saved cell ends up with no explicit style - matching the other untouched cells.
$numericLookingCells = @{
    'D4' = '1.00'
    'D5' = '597.58'
    'D6' = '149.84'
    'D7' = '1.00'
    'D9' = '0.519'
    'D12' = '0.458'
    'D13' = '0.0000233'
    'D14' = '34.56'
    'D18' = '7.01'
    'D20' = '448.63'
    'D21' = '14.19'
    'D22' = '0.691'
    'D23' = '7.45'
    'D24' = '82.35'
    'D25' = '10.86'
    'D30' = '7.29'
    'D31' = '1.00'
    'D32' = '2.15'
    'D33' = '27.55'
    'D37' = '5.85'
    'D38' = '3.03'
    'D39' = '2.07'
    'D40' = '50.12'
    'D44' = '391.32'
    'D45' = '40.10'
    'D46' = '0.0355'
    'D48' = '133.06'
}
foreach ($addr in $numericLookingCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $numericLookingCells[$addr]
    $rng.Style = "Normal"
}

# Remaining cells (coin names, links, percent strings, multi-dot prices, etc.) are
# not parsed as numbers by Excel, so a plain Value assignment is sufficient.
$textCells = @{
    'D2' = '62.594.73'
    'E2' = '  +1.77%  '
    'D3' = '3.021.54'
    'E3' = '  +2.26%  '
    'E4' = '  +0.10%  '
    'E5' = '  +2.17%  '
    'E6' = '  +6.51%  '
    'E7' = '  +0.05%  '
    'D8' = '3.019.70'
    'E8' = '  +2.20%  '
    'E9' = '  -0.08%  '
    'E10' = '  +12.17%  '
    'E11' = '  +4.80%  '
    'E12' = '  +0.72%  '
    'E13' = '  +4.07%  '
    'E14' = '  +2.39%  '
    'E15' = '  +2.71%  '
    'D16' = '3.524.98'
    'E16' = '  +2.44%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '62.626.52'
    'E17' = '  +1.87%  '
    'B18' = 'Polkadot'
    'C18' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'E18' = '  +0.71%  '
    'D19' = '3.023.38'
    'E19' = '  +2.65%  '
    'E20' = '  +0.23%  '
    'E21' = '  +2.93%  '
    'E22' = '  +1.86%  '
    'E23' = '  +2.52%  '
    'E24' = '  +1.99%  '
    'E25' = '  +13.74%  '
    'E26' = '  +5.37%  '
    'E29' = '  +4.15%  '
    'B30' = 'NEARProtocol'
    'C30' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E30' = '  +6.60%  '
    'B31' = 'FirstDigitalUSD'
    'C31' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'E31' = '  +0.20%  '
    'E32' = '  +5.05%  '
    'E33' = '  +2.04%  '
    'E34' = '  +3.47%  '
    'D35' = '0.0₃0853'
    'E35' = '  +11.66%  '
    'E36' = '  +2.31%  '
    'E37' = '  +3.59%  '
    'E38' = '  +9.54%  '
    'E39' = '  +0.52%  '
    'E40' = '  +0.25%  '
    'E41' = '  -0.74%  '
    'E42' = '  +4.64%  '
    'E43' = '  +9.28%  '
    'E44' = '  +1.67%  '
    'E45' = '  +9.05%  '
    'E46' = '  +1.01%  '
    'D47' = '2.742.21'
    'E47' = '  +1.33%  '
    'E48' = '  +2.65%  '
    'E49' = '  +0.10%  '
    'E50' = '  +2.15%  '
    'E51' = '  +0.43%  '
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
